$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(84210, 95730, 103400, 87450, 102780, 98720, 105360, 92670, 87690, 94420, 101530, 91200)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("A1:B13").Select()
